$wb = $excel.ActiveWorkbook

# Rename first sheet to match usage in script
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "WCC All Approved Mentors"

# Make the first sheet the active/selected sheet (instead of the second one)
$ws1.Activate()
